# Update the data-type column ("形態") for the CreateDate / LastUpdate rows
# on the "DBD" layout sheet: DATE -> TIMESTAMP.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("D12").Value = "TIMESTAMP"
$ws.Range("D14").Value = "TIMESTAMP"
